# Updated CHE_grids model - 2025-08-19 00:07
# Re-shuffle of the "grid_cell" (AG) column on the "solar" sheet: each row's
# process (AC/AD/AE) stays the same, but the grid-cell id label in AG is
# reassigned to a new CHE_NN value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("solar")

$newGridCells = @{
    4  = "CHE_3"
    5  = "CHE_0"
    6  = "CHE_2"
    7  = "CHE_17"
    8  = "CHE_19"
    9  = "CHE_23"
    10 = "CHE_10"
    11 = "CHE_22"
    12 = "CHE_24"
    13 = "CHE_8"
    14 = "CHE_5"
    15 = "CHE_11"
    16 = "CHE_15"
    17 = "CHE_25"
    18 = "CHE_14"
    19 = "CHE_18"
    20 = "CHE_7"
    21 = "CHE_9"
    22 = "CHE_21"
    23 = "CHE_4"
    24 = "CHE_12"
    25 = "CHE_20"
    26 = "CHE_1"
    27 = "CHE_6"
    28 = "CHE_13"
}

foreach ($row in $newGridCells.Keys) {
    $ws.Range("AG" + $row).Value = $newGridCells[$row]
}
